$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table 3 ("Caso de Uso" / Peso / N Casos de Uso / Resultado)
#    Medio row: N Casos de Uso 0 -> 1 ; Resultado 0 -> 10
#    Totals row: Resultado 20 -> 30
# ---------------------------------------------------------------------------
$tUseCases = $d.Tables.Item(3)
$tUseCases.Cell(3, 3).Range.Text = "1"
$tUseCases.Cell(3, 4).Range.Text = "10"
$tUseCases.Cell(5, 4).Range.Text = "30"

# ---------------------------------------------------------------------------
# 2) "UUCP = 20 + 6 = 26"  ->  "UUCP = 30 + 6 = 36"
# ---------------------------------------------------------------------------
$pUUCP = $d.Paragraphs.Item(197).Range
$pUUCP.Find.Execute("20 + 6 = 26", $true, $false, $false, $false, $false, $true, 1, $false, "30 + 6 = 36", 2)

# ---------------------------------------------------------------------------
# 3) "UCP = 26 * 0.87 * 0.995 = 22.5069" -> "UCP = 36 * 0,87 * 0,995 = 31,1634"
# ---------------------------------------------------------------------------
$pUCP = $d.Paragraphs.Item(198).Range
$pUCP.Find.Execute("26 * 0.87 * 0.995 = 22.5069", $true, $false, $false, $false, $false, $true, 1, $false, "36 * 0,87 * 0,995 = 31,1634", 2)

# ---------------------------------------------------------------------------
# 4) "E = UCP * CF =  22.5069* 20 = 450.138 [Horas - Hombre]"
#    -> "E = UCP * CF =  31,1634* 20 = 623,268  [Horas - Hombre]"
# ---------------------------------------------------------------------------
$pE = $d.Paragraphs.Item(204).Range
$pE.Find.Execute("22.5069* 20 = 450.138 ", $true, $false, $false, $false, $false, $true, 1, $false, "31,1634* 20 = 623,268  ", 2)

# ---------------------------------------------------------------------------
# 5) "TDEV = 450.138   [Horas - Hombre]   /    3 [Hombre]  = 150.046 Horas..."
#    -> "TDEV = 623,268   [Horas - Hombre]   /    3 [Hombre]  = 207,756 Horas..."
# ---------------------------------------------------------------------------
$pTDEV = $d.Paragraphs.Item(208).Range
$pTDEV.Find.Execute("450.138", $true, $false, $false, $false, $false, $true, 1, $false, "623,268", 2)
$pTDEV2 = $d.Paragraphs.Item(208).Range
$pTDEV2.Find.Execute("150.046", $true, $false, $false, $false, $false, $true, 1, $false, "207,756", 2)

# ---------------------------------------------------------------------------
# 6) "150.046 / 7 = 21,44  semanas ...  5,4 meses."
#    -> "207,756 / 7 = 29,68  semanas ...  7,42 meses."
# ---------------------------------------------------------------------------
$pWeeks = $d.Paragraphs.Item(211).Range
$pWeeks.Find.Execute("150.046", $true, $false, $false, $false, $false, $true, 1, $false, "207,756", 2)
$pWeeks2 = $d.Paragraphs.Item(211).Range
$pWeeks2.Find.Execute("21,44", $true, $false, $false, $false, $false, $true, 1, $false, "29,68", 2)
$pWeeks3 = $d.Paragraphs.Item(211).Range
$pWeeks3.Find.Execute("5,4 ", $true, $false, $false, $false, $false, $true, 1, $false, "7,42 ", 2)

# ---------------------------------------------------------------------------
# 7) Table 5 (Actividad / Porcentaje)
# ---------------------------------------------------------------------------
$tAct = $d.Tables.Item(5)
$tAct.Cell(2, 2).Range.Text = "10.00% = 2,9 weeks"
$tAct.Cell(3, 2).Range.Text = "20.00% = 5,9 weeks"
$tAct.Cell(4, 2).Range.Text = "40.00% = 11,9 +2 (Vacaciones)weeks"
$tAct.Cell(5, 2).Range.Text = "15.00% = 4,5 weeks"
$tAct.Cell(6, 2).Range.Text = "15.00% = 4,5 weeks"

# ---------------------------------------------------------------------------
# 8) Move the hidden "_GoBack" bookmark from inside the activities table to
#    the very end of the document (best effort; ignore failures since this
#    bookmark is an internal Word navigation marker, not visible content).
# ---------------------------------------------------------------------------
try {
    $endPos = $d.Content.End - 1
    $endRng = $d.Range($endPos, $endPos)
    $d.Bookmarks.Add("_GoBack", $endRng)
} catch {
    Write-Host "Bookmark move skipped: $_"
}
